# Add a new model run log entry "2023_TM160_IPA_28" to the ModelRuns sheet.
# This inserts a new row 40 (pushing the existing rows 40-47 down to 41-48),
# copying the format of the row above (row 39, the prior "current" run) and
# filling in the new run's data, then updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 40; Excel copies formatting from the row above (row 39)
$ws.Rows.Item(40).Insert()

# Most fields repeat the same values as the row above (row 39), since this
# new run continues the same category/run_set/model_machine/etc.
$ws.Cells.Item(40, 1).Value = $ws.Cells.Item(39, 1).Value2     # year -> 2023
$ws.Cells.Item(40, 2).Value = "2023_TM160_IPA_28"               # directory
$ws.Cells.Item(40, 3).Value = $ws.Cells.Item(39, 3).Value2      # run_set
$ws.Cells.Item(40, 4).Value = $ws.Cells.Item(39, 4).Value2      # category
$ws.Cells.Item(40, 5).Value = "Added BART hesistancy to Caltrain and ferry, WFH remains at ~25%"  # description
$ws.Cells.Item(40, 6).Value = $ws.Cells.Item(39, 6).Value2      # urbansim_runid
$ws.Cells.Item(40, 7).Value = $ws.Cells.Item(39, 7).Value2      # status
$ws.Cells.Item(40, 8).Value = $ws.Cells.Item(39, 8).Value2      # network
$ws.Cells.Item(40, 9).Value = $ws.Cells.Item(39, 9).Value2      # landuse_path
$ws.Cells.Item(40, 10).Value = $ws.Cells.Item(39, 10).Value2    # model_machine
$ws.Cells.Item(40, 11).Value = $ws.Cells.Item(39, 11).Value2    # Asana Link
$ws.Cells.Item(40, 12).Value = 17.77
$ws.Cells.Item(40, 13).Value = $ws.Cells.Item(39, 13).Value2
$ws.Cells.Item(40, 14).Value = $ws.Cells.Item(39, 14).Value2
$ws.Cells.Item(40, 15).Value = 0.94
$ws.Cells.Item(40, 16).Value = 0.855
$ws.Cells.Item(40, 17).Value = 120
$ws.Cells.Item(40, 18).Value = 0
$ws.Cells.Item(40, 19).Value = 45

# Update the current selection to the new row's run-name cell
$ws.Range("B40").Select()
